# Applies the edits described by the commit "v0.1 - something works but some issues"
# to the steel_data workbook.
#
# Summary of content changes:
#  - baseline: column F header changes from "remained_lifespan" to "introduced_year",
#    and its values are replaced with the actual year each technology was introduced.
#  - capex: ESF (row 3) future capex (2039-2050 / cols P:AA) set to 0, and
#    Electricity (row 4) future capex (2039-2050 / cols P:AA) set to 20.
#    (opex/renewal recompute automatically since they reference capex via formulas.)
#  - fuel_cost: Electricity (row 4) fuel cost for every year (cols B:AA) set to 0.
#  - Selection/active-cell bookkeeping on a few sheets (capex, technology, fuel_cost)
#    changed as a side effect of the user's editing session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# baseline: rename/repurpose the last column and give it real content
# ---------------------------------------------------------------------------
$baseline = $wb.Worksheets.Item("baseline")
$baseline.Range("F1").Value2 = "introduced_year"
$baseline.Range("F2").Value2 = 2015
$baseline.Range("F3").Value2 = 2020
$baseline.Range("F4").Value2 = 2023
$baseline.Range("F5").Value2 = 2021

# ---------------------------------------------------------------------------
# capex: zero out ESF's future-year capex, flatten Electricity's future-year capex to 20
# ---------------------------------------------------------------------------
$capex = $wb.Worksheets.Item("capex")
$capex.Range("P3:AA3").Value2 = 0
$capex.Range("P4:AA4").Value2 = 20
$capex.Range("D8").Select()

# ---------------------------------------------------------------------------
# technology: move the active selection
# ---------------------------------------------------------------------------
$technology = $wb.Worksheets.Item("technology")
$technology.Range("C3").Select()

# ---------------------------------------------------------------------------
# fuel_cost: zero out Electricity's fuel cost across every year
# ---------------------------------------------------------------------------
$fuelCost = $wb.Worksheets.Item("fuel_cost")
$fuelCost.Range("B4:AA4").Value2 = 0
$fuelCost.Range("H9").Select()

# ---------------------------------------------------------------------------
# restore "baseline" as the active sheet/tab without disturbing its own
# remembered selection (F6)
# ---------------------------------------------------------------------------
$baseline.Select()
